$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @("A","B","C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T")

$data = @(
    @{row=2; vals=@('ECs', 'Gal', 'Galr2', 'ECs', 3, 1, 1.776285, 5.328855000000001, 0.2250740306326953, 0.2250740306326953, 2, 0.6666666666666666, 0.402287, 1.206861, 0.09937001049721769, 0.09937001049721769, 0.7145763637950001, 6.431187274155, 0.02236560878662203, 0.02236560878662203)}
    @{row=3; vals=@('ECs', 'Gal', 'Galr2', 'FAPs', 3, 1, 1.776285, 5.328855000000001, 0.2250740306326953, 0.2250740306326953, 3, 1, 1.920356666666667, 5.76107, 0.4743525446386999, 0.4743525446386999, 3.411100741650001, 30.69990667485001, 0.1067644391627077, 0.1067644391627077)}
    @{row=4; vals=@('ECs', 'Gal', 'Galr2', 'MuSCs', 3, 1, 1.776285, 5.328855000000001, 0.2250740306326953, 0.2250740306326953, 3, 1, 1.553186, 4.659558000000001, 0.3836567155392334, 0.3836567155392334, 2.758900994010001, 24.83010894609001, 0.08635116334571669, 0.08635116334571669)}
    @{row=5; vals=@('ECs', 'Gal', 'Galr2', 'Resolving-Mac', 3, 1, 1.776285, 5.328855000000001, 0.2250740306326953, 0.2250740306326953, 2, 0.6666666666666666, 0.1725446666666667, 0.517634, 0.04262072932484916, 0.04262072932484916, 0.3064885032300001, 2.758396529070001, 0.009592819337648916, 0.009592819337648916)}
    @{row=6; vals=@('FAPs', 'Gal', 'Galr2', 'ECs', 3, 1, 1.407624, 4.222872, 0.1783607964348723, 0.1783607964348722, 2, 0.6666666666666666, 0.402287, 1.206861, 0.09937001049721769, 0.09937001049721769, 0.566268836088, 5.096419524791999, 0.01772371421402536, 0.01772371421402536)}
    @{row=7; vals=@('FAPs', 'Gal', 'Galr2', 'FAPs', 3, 1, 1.407624, 4.222872, 0.1783607964348723, 0.1783607964348722, 3, 1, 1.920356666666667, 5.76107, 0.4743525446386999, 0.4743525446386999, 2.70314013256, 24.32826119304, 0.08460589765266681, 0.0846058976526668)}
    @{row=8; vals=@('FAPs', 'Gal', 'Galr2', 'MuSCs', 3, 1, 1.407624, 4.222872, 0.1783607964348723, 0.1783607964348722, 3, 1, 1.553186, 4.659558000000001, 0.3836567155392334, 0.3836567155392334, 2.186301890064, 19.676717010576, 0.06842931734116489, 0.06842931734116489)}
    @{row=9; vals=@('FAPs', 'Gal', 'Galr2', 'Resolving-Mac', 3, 1, 1.407624, 4.222872, 0.1783607964348723, 0.1783607964348722, 2, 0.6666666666666666, 0.1725446666666667, 0.517634, 0.04262072932484916, 0.04262072932484916, 0.242878013872, 2.185902124848, 0.007601867227015213, 0.007601867227015211)}
    @{row=10; vals=@('MuSCs', 'Gal', 'Galr2', 'ECs', 3, 1, 4.248107333333333, 12.744322, 0.5382799720054182, 0.5382799720054181, 2, 0.6666666666666666, 0.402287, 1.206861, 0.09937001049721769, 0.09937001049721769, 1.708958354804667, 15.380625193242, 0.05348888646862045, 0.05348888646862044)}
    @{row=11; vals=@('MuSCs', 'Gal', 'Galr2', 'FAPs', 3, 1, 4.248107333333333, 12.744322, 0.5382799720054182, 0.5382799720054181, 3, 1, 1.920356666666667, 5.76107, 0.4743525446386999, 0.4743525446386999, 8.157881238282222, 73.42093114454, 0.2553344744488182, 0.2553344744488182)}
    @{row=12; vals=@('MuSCs', 'Gal', 'Galr2', 'MuSCs', 3, 1, 4.248107333333333, 12.744322, 0.5382799720054182, 0.5382799720054181, 3, 1, 1.553186, 4.659558000000001, 0.3836567155392334, 0.3836567155392334, 6.598100836630667, 59.38290752967601, 0.2065147261001492, 0.2065147261001492)}
    @{row=13; vals=@('MuSCs', 'Gal', 'Galr2', 'Resolving-Mac', 3, 1, 4.248107333333333, 12.744322, 0.5382799720054182, 0.5382799720054181, 2, 0.6666666666666666, 0.1725446666666667, 0.517634, 0.04262072932484916, 0.04262072932484916, 0.7329882637942223, 6.596894374148, 0.02294188498783032, 0.02294188498783031)}
    @{row=14; vals=@('Resolving-Mac', 'Gal', 'Galr2', 'ECs', 2, 0.6666666666666666, 0.459987, 1.379961, 0.05828520092701431, 0.0582852009270143, 2, 0.6666666666666666, 0.402287, 1.206861, 0.09937001049721769, 0.09937001049721769, 0.185046790269, 1.665421112421, 0.005791801027949854, 0.005791801027949853)}
    @{row=15; vals=@('Resolving-Mac', 'Gal', 'Galr2', 'FAPs', 2, 0.6666666666666666, 0.459987, 1.379961, 0.05828520092701431, 0.0582852009270143, 3, 1, 1.920356666666667, 5.76107, 0.4743525446386999, 0.4743525446386999, 0.88333910203, 7.95005191827, 0.02764773337450714, 0.02764773337450714)}
    @{row=16; vals=@('Resolving-Mac', 'Gal', 'Galr2', 'MuSCs', 2, 0.6666666666666666, 0.459987, 1.379961, 0.05828520092701431, 0.0582852009270143, 3, 1, 1.553186, 4.659558000000001, 0.3836567155392334, 0.3836567155392334, 0.714445368582, 6.430008317238, 0.02236150875220259, 0.02236150875220259)}
    @{row=17; vals=@('Resolving-Mac', 'Gal', 'Galr2', 'Resolving-Mac', 2, 0.6666666666666666, 0.459987, 1.379961, 0.05828520092701431, 0.0582852009270143, 2, 0.6666666666666666, 0.1725446666666667, 0.517634, 0.04262072932484916, 0.04262072932484916, 0.079368303586, 0.714314732274, 0.002484157772354724, 0.002484157772354724)}
)

foreach ($item in $data) {
    $r = $item.row
    $vals = $item.vals
    for ($i = 0; $i -lt $cols.Length; $i++) {
        $ws.Range($cols[$i] + $r).Value = $vals[$i]
    }
}
